$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9975859522819519
$ws.Range("B1").Value = 1.817065834999084
$ws.Range("C1").Value = 6.884543418884277
$ws.Range("D1").Value = 2.89591646194458
$ws.Range("E1").Value = 0.4142286777496338
